# Insert a new data row at row 314 (shifting rows 314:387 down to 315:388)
# and populate it with the new "Perejil" price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 314, shifting existing rows down.
$ws.Rows.Item(314).Insert(-4121)   # -4121 == xlShiftDown

# Fill in the values for the newly inserted row 314.
$ws.Cells.Item(314, 1).Value  = 9
$ws.Cells.Item(314, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(314, 3).Value  = "Metropolitana"
$ws.Cells.Item(314, 4).Value  = 44722
$ws.Cells.Item(314, 5).Value  = 13
$ws.Cells.Item(314, 6).Value  = 100112044
$ws.Cells.Item(314, 7).Value  = "Perejil"
$ws.Cells.Item(314, 8).Value  = "Sin especificar"
$ws.Cells.Item(314, 9).Value  = "Primera"
$ws.Cells.Item(314, 10).Value = 61
$ws.Cells.Item(314, 11).Value = 8000
$ws.Cells.Item(314, 12).Value = 9000
$ws.Cells.Item(314, 13).Value = 8508
$ws.Cells.Item(314, 14).Value = "$/docena de atados"
$ws.Cells.Item(314, 15).Value = "Región Metropolitana"
$ws.Cells.Item(314, 16).Value = 2836
$ws.Cells.Item(314, 17).Value = 3
$ws.Cells.Item(314, 18).Value = "Hortaliza"
